$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false
$ws = $wb.Worksheets("Login")

# Remove the extra sheets ("User" and "Program") - only "Login" remains.
$wb.Worksheets("User").Delete()
$wb.Worksheets("Program").Delete()

# Remove the existing hyperlink on A2 (its text stays, but the link itself moves later).
$ws.Hyperlinks.Delete()

# Insert a new row before row 3 so the second email/password pair has its own row.
$ws.Rows("3").Insert()

# Row 2: email (A2) / original password (B2)
$ws.Range("A2").Value = "sdetorganizers@gmail.com"
$ws.Range("B2").Value = "UIHackathon@02"

# Row 3: email again (A3) / new password with hyperlink (B3)
$ws.Range("A3").Value = "sdetorganizers@gmail.com"
$ws.Range("B3").Value = "UIHackathon@021231"

$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:sdetorganizers@gmail.com", "", "", "UIHackathon@021231")

# Adjust column B width and the selected cell, matching the recorded view state.
$ws.Columns("B").ColumnWidth = 17.44140625
$ws.Range("B8").Select()

$wb.Save()
